$d = $word.ActiveDocument

# --- Step 1: paragraph 1 ("This is a Microsoft word document.") ---
# Append two trailing spaces to the existing run's text, then append three
# new runs colored C00000: "(This is a change – Ve", "rsion for branch
# alternate", ")"
$p1 = $d.Paragraphs.First
$p1End = $p1.Range.End - 1          # position right before paragraph mark

$spacesRange = $d.Range($p1End, $p1End)
$spacesRange.InsertAfter("  ")

$pos = $p1End + 2

$seg1 = "(This is a change – Ve"
$r1 = $d.Range($pos, $pos)
$r1.InsertAfter($seg1)
$seg1Range = $d.Range($pos, $pos + $seg1.Length)
$seg1Range.Font.Color = 192
$pos = $pos + $seg1.Length

$seg2 = "rsion for branch alternate"
$r2 = $d.Range($pos, $pos)
$r2.InsertAfter($seg2)
$seg2Range = $d.Range($pos, $pos + $seg2.Length)
$seg2Range.Font.Color = 192
$pos = $pos + $seg2.Length

$seg3 = ")"
$r3 = $d.Range($pos, $pos)
$r3.InsertAfter($seg3)
$seg3Range = $d.Range($pos, $pos + $seg3.Length)
$seg3Range.Font.Color = 192
$pos = $pos + $seg3.Length

# --- Step 2: new empty (shaded) paragraph right after
#     "It will be treated as a binary file by Git." ---
$p2 = $d.Paragraphs(2)
$p2End = $p2.Range.End - 1          # position right before paragraph mark

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p>'
$paraInsertRange = $d.Range($p2End, $p2End)
$paraInsertRange.InsertXML($newParaXml)

Write-Output "done"
